# Update the NMDP ABO codes ValueSet workbook for the 0.1.1 release:
#  - rename the "Include from NMDP Blood Group" sheet to "Include #0"
#  - bump Version 0.1.0 -> 0.1.1
#  - update the Date field
#  - add a "Jurisdiction" metadata row (standard FHIR IG metadata field,
#    inserted right after "Contact" and before "Description")

$wb = $excel.ActiveWorkbook

# 1. Rename the include sheet.
$wsInclude = $wb.Worksheets.Item("Include from NMDP Blood Group")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet.
$ws = $wb.Worksheets.Item("Metadata")

# Version (row 3) : 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Date (row 8)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10),
# before "Description" (old row 11) - this pushes Description / Purpose /
# Copyright / Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Carry the body-row formatting (border/alignment) down onto the new row
# by copying the format from the row above it.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"

Write-Host "Workbook updated for 0.1.1 release"
